# [PHOENIX-5910] added search trade feature
#
# Populates the new "searchTradeDeatils" sheet (sheet6) with the lookup
# rows used by the search-trade functional tests, updates the sample
# "licenseClosure" row's trade category, and flips which sheet/tab is
# active/selected in the workbook.

$wb = $excel.ActiveWorkbook

$wsSearch = $wb.Worksheets.Item("searchTradeDeatils")
$wsClosure = $wb.Worksheets.Item("licenseClosure")

# --- searchTradeDeatils (sheet6): add the dataRow / searchValue table ---
$wsSearch.Range("B1").Value = "searchValue"
$wsSearch.Range("A1").Value = "dataRow"
$wsSearch.Range("A2").Value = "searchWithApplicationNumber"
$wsSearch.Range("A3").Value = "searchWithLicenseNumber"
$wsSearch.Range("B2").Value = "00393-2017-PG "
$wsSearch.Range("B3").Value = "TL/08360/2016"

$wsSearch.Columns.Item(1).AutoFit() | Out-Null
$wsSearch.Columns.Item(2).AutoFit() | Out-Null

$wsSearch.Range("C9").Select() | Out-Null

# --- licenseClosure (sheet5): change the sample trade category ---
$wsClosure.Range("C2").Value = "Veterinary Trades"

$wsClosure.Range("D11").Select() | Out-Null

# licenseClosure becomes the active/selected sheet/tab
$wsClosure.Activate() | Out-Null

Write-Host "done"
